$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "304.37" },
    @{ Cell = "E2"; Value = "5.76%" },
    @{ Cell = "D3"; Value = "35.01" },
    @{ Cell = "E3"; Value = "12.81%" },
    @{ Cell = "D4"; Value = "5.169" },
    @{ Cell = "E4"; Value = "5.05%" },
    @{ Cell = "D5"; Value = "0.07770" },
    @{ Cell = "E5"; Value = "6.33%" },
    @{ Cell = "D6"; Value = "2.297" },
    @{ Cell = "E6"; Value = "1.91%" },
    @{ Cell = "D8"; Value = "4.005" },
    @{ Cell = "E8"; Value = "7.66%" },
    @{ Cell = "D9"; Value = "0.9296" },
    @{ Cell = "E9"; Value = "2.83%" },
    @{ Cell = "D10"; Value = "0.1013" },
    @{ Cell = "E10"; Value = "10.92%" },
    @{ Cell = "D11"; Value = "0.1824" },
    @{ Cell = "E11"; Value = "8.08%" },
    @{ Cell = "D12"; Value = "0.08666" },
    @{ Cell = "E12"; Value = "5.81%" },
    @{ Cell = "D13"; Value = "0.03443" },
    @{ Cell = "E13"; Value = "10.31%" },
    @{ Cell = "D14"; Value = "0.09880" },
    @{ Cell = "E14"; Value = "-0.61%" },
    @{ Cell = "D15"; Value = "0.001482" },
    @{ Cell = "E15"; Value = "-1.45%" },
    @{ Cell = "D16"; Value = "0.04614" },
    @{ Cell = "E16"; Value = "2.26%" },
    @{ Cell = "D17"; Value = "0.005852" },
    @{ Cell = "E17"; Value = "2.57%" },
    @{ Cell = "D18"; Value = "3.506" },
    @{ Cell = "E18"; Value = "-0.18%" },
    @{ Cell = "D19"; Value = "2.112" },
    @{ Cell = "E19"; Value = "3.14%" },
    @{ Cell = "D20"; Value = "0.3419" },
    @{ Cell = "E20"; Value = "2.64%" },
    @{ Cell = "D21"; Value = "0.1328" },
    @{ Cell = "E21"; Value = "2.24%" },
    @{ Cell = "D22"; Value = "4.643" },
    @{ Cell = "E22"; Value = "10.29%" },
    @{ Cell = "D23"; Value = "0.2341" },
    @{ Cell = "E23"; Value = "11.44%" },
    @{ Cell = "D24"; Value = "0.001225" },
    @{ Cell = "E24"; Value = "1.08%" },
    @{ Cell = "D25"; Value = "0.004418" },
    @{ Cell = "E25"; Value = "6.34%" },
    @{ Cell = "D26"; Value = "0.0001304" },
    @{ Cell = "E26"; Value = "0.26%" },
    @{ Cell = "D27"; Value = "0.0003419" },
    @{ Cell = "E27"; Value = "0.67%" },
    @{ Cell = "D39"; Value = "0.01765" },
    @{ Cell = "E39"; Value = "12.33%" },
    @{ Cell = "D40"; Value = "0.04731" },
    @{ Cell = "E40"; Value = "6.54%" },
    @{ Cell = "D41"; Value = "0.007657" },
    @{ Cell = "E41"; Value = "4.62%" },
    @{ Cell = "D42"; Value = "0.1406" },
    @{ Cell = "E42"; Value = "6.13%" },
    @{ Cell = "D43"; Value = "0.007079" },
    @{ Cell = "E43"; Value = "-25.72%" },
    @{ Cell = "D44"; Value = "0.002179" },
    @{ Cell = "E44"; Value = "-1.87%" },
    @{ Cell = "D45"; Value = "0.009215" },
    @{ Cell = "E45"; Value = "2.99%" },
    @{ Cell = "D46"; Value = "0.00005920" },
    @{ Cell = "E46"; Value = "-3.09%" },
    @{ Cell = "E47"; Value = "0.18%" },
    @{ Cell = "E48"; Value = "19.10%" },
    @{ Cell = "D49"; Value = "0.002704" },
    @{ Cell = "E49"; Value = "35.10%" },
    @{ Cell = "D50"; Value = "0.00002104" },
    @{ Cell = "E50"; Value = "0.18%" },
    @{ Cell = "D51"; Value = "0.0002004" },
    @{ Cell = "E51"; Value = "0.18%" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
